$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1. sheet1 new cells (order matters for shared-string allocation: URL must be the first new string)
$ws1.Range("K2").Value = "https://www.facebook.com/profile.php?id=100013637047044"
$ws1.Range("I2").Value = $true

# 2. sheet2: insert column before H (after G/"ВУЗ")
$ws2.Columns.Item(8).Insert()
$ws2.Range("H1").Value = "Город(ВУЗ)"
$ws2.Range("H2").Value = "Москва"
$ws2.Range("H3").Value = "Москва"

# 3. sheet2: insert column before L (after K/"Школа", which is now at column K)
$ws2.Columns.Item(12).Insert()
$ws2.Range("L1").Value = "Город(школа)"
$ws2.Range("L2").Value = "Москва"
$ws2.Range("L3").Value = "Москва"

# 4. sheet2: fix K2 (now holding stale "СШ№2") to new school name
$ws2.Range("K2").Value = "сш №209"

# 5. activate sheet2 and set selection
$ws2.Activate()
$ws2.Range("K9").Select()
